$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.447252
$ws.Range("H2").Value = 61.341756
$ws.Range("I2").Value = 0.8699145605694745
$ws.Range("J2").Value = 0.8770588936480435
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.014513
$ws.Range("N2").Value = 0.043539
$ws.Range("O2").Value = 0.006538124153480057
$ws.Range("P2").Value = 0.009244530017173054
$ws.Range("Q2").Value = 0.296750968276
$ws.Range("R2").Value = 2.670758714484
$ws.Range("S2").Value = 0.005687609399923271
$ws.Range("T2").Value = 0.008107997269157927

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.447252
$ws.Range("H3").Value = 61.341756
$ws.Range("I3").Value = 0.8699145605694745
$ws.Range("J3").Value = 0.8770588936480435
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2421403333333333
$ws.Range("N3").Value = 0.726421
$ws.Range("O3").Value = 0.1090845147039467
$ws.Range("P3").Value = 0.154239204841748
$ws.Range("Q3").Value = 4.951104415030666
$ws.Range("R3").Value = 44.559939735276
$ws.Range("S3").Value = 0.0948942076736182
$ws.Range("T3").Value = 0.1352768663556574

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha7"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.447252
$ws.Range("H4").Value = 61.341756
$ws.Range("I4").Value = 0.8699145605694745
$ws.Range("J4").Value = 0.8770588936480435
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.013551
$ws.Range("N4").Value = 0.040653
$ws.Range("O4").Value = 0.006104741983312083
$ws.Range("P4").Value = 0.008631752653669954
$ws.Range("Q4").Value = 0.277080711852
$ws.Range("R4").Value = 2.493726406668
$ws.Range("S4").Value = 0.005310603939802953
$ws.Range("T4").Value = 0.007570555432671333

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha7"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.447252
$ws.Range("H5").Value = 61.341756
$ws.Range("I5").Value = 0.8699145605694745
$ws.Range("J5").Value = 0.8770588936480435
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.9495455
$ws.Range("N5").Value = 3.899091
$ws.Range("O5").Value = 0.878272619159261
$ws.Range("P5").Value = 0.827884512487409
$ws.Range("Q5").Value = 39.862848123966
$ws.Range("R5").Value = 239.177088743796
$ws.Range("S5").Value = 0.76402213955613
$ws.Range("T5").Value = 0.7261034745905568

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha7"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.483247333333333
$ws.Range("H6").Value = 7.449742
$ws.Range("I6").Value = 0.105648084777455
$ws.Range("J6").Value = 0.1065157390747562
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.014513
$ws.Range("N6").Value = 0.043539
$ws.Range("O6").Value = 0.006538124153480057
$ws.Range("P6").Value = 0.009244530017173054
$ws.Range("Q6").Value = 0.03603936854866667
$ws.Range("R6").Value = 0.324354316938
$ws.Range("S6").Value = 0.0006907402948523872
$ws.Range("T6").Value = 0.0009846879471779569

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha7"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.483247333333333
$ws.Range("H7").Value = 7.449742
$ws.Range("I7").Value = 0.105648084777455
$ws.Range("J7").Value = 0.1065157390747562
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2421403333333333
$ws.Range("N7").Value = 0.726421
$ws.Range("O7").Value = 0.1090845147039467
$ws.Range("P7").Value = 0.154239204841748
$ws.Range("Q7").Value = 0.6012943370424444
$ws.Range("R7").Value = 5.411649033382
$ws.Range("S7").Value = 0.0115245700573501
$ws.Range("T7").Value = 0.01642890289802151

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha7"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.483247333333333
$ws.Range("H8").Value = 7.449742
$ws.Range("I8").Value = 0.105648084777455
$ws.Range("J8").Value = 0.1065157390747562
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.013551
$ws.Range("N8").Value = 0.040653
$ws.Range("O8").Value = 0.006104741983312083
$ws.Range("P8").Value = 0.008631752653669954
$ws.Range("Q8").Value = 0.033650484614
$ws.Range("R8").Value = 0.302854361526
$ws.Range("S8").Value = 0.0006449542985974436
$ws.Range("T8").Value = 0.0009194175134161436

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha7"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.483247333333333
$ws.Range("H9").Value = 7.449742
$ws.Range("I9").Value = 0.105648084777455
$ws.Range("J9").Value = 0.1065157390747562
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.9495455
$ws.Range("N9").Value = 3.899091
$ws.Range("O9").Value = 0.878272619159261
$ws.Range("P9").Value = 0.827884512487409
$ws.Range("Q9").Value = 4.841203664087001
$ws.Range("R9").Value = 29.047221984522
$ws.Range("S9").Value = 0.09278782012665505
$ws.Range("T9").Value = 0.08818273071614063

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha7"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5743975
$ws.Range("H10").Value = 1.148795
$ws.Range("I10").Value = 0.02443735465307048
$ws.Range("J10").Value = 0.01642536727720028
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.014513
$ws.Range("N10").Value = 0.043539
$ws.Range("O10").Value = 0.006538124153480057
$ws.Range("P10").Value = 0.009244530017173054
$ws.Range("Q10").Value = 0.0083362309175
$ws.Range("R10").Value = 0.050017385505
$ws.Range("S10").Value = 0.0001597744587043984
$ws.Range("T10").Value = 0.0001518448008371701

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Efna1"
$ws.Range("C11").Value = "Epha7"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5743975
$ws.Range("H11").Value = 1.148795
$ws.Range("I11").Value = 0.02443735465307048
$ws.Range("J11").Value = 0.01642536727720028
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2421403333333333
$ws.Range("N11").Value = 0.726421
$ws.Range("O11").Value = 0.1090845147039467
$ws.Range("P11").Value = 0.154239204841748
$ws.Range("Q11").Value = 0.1390848021158333
$ws.Range("R11").Value = 0.834508812695
$ws.Range("S11").Value = 0.002665736972978428
$ws.Range("T11").Value = 0.002533435588069039

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Efna1"
$ws.Range("C12").Value = "Epha7"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5743975
$ws.Range("H12").Value = 1.148795
$ws.Range("I12").Value = 0.02443735465307048
$ws.Range("J12").Value = 0.01642536727720028
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.013551
$ws.Range("N12").Value = 0.040653
$ws.Range("O12").Value = 0.006104741983312083
$ws.Range("P12").Value = 0.008631752653669954
$ws.Range("Q12").Value = 0.0077836605225
$ws.Range("R12").Value = 0.046701963135
$ws.Range("S12").Value = 0.0001491837449116863
$ws.Range("T12").Value = 0.0001417797075824772

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Efna1"
$ws.Range("C13").Value = "Epha7"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5743975
$ws.Range("H13").Value = 1.148795
$ws.Range("I13").Value = 0.02443735465307048
$ws.Range("J13").Value = 0.01642536727720028
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.9495455
$ws.Range("N13").Value = 3.899091
$ws.Range("O13").Value = 0.878272619159261
$ws.Range("P13").Value = 0.827884512487409
$ws.Range("Q13").Value = 1.11981406133625
$ws.Range("R13").Value = 4.479256245345001
$ws.Range("S13").Value = 0.02146265947647597
$ws.Range("T13").Value = 0.0135983071807116
